$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values for rows 2-10 from 45175 (2023-09-06)
# to 45183 (2023-09-14), keeping existing cell formatting.
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
